# Auto-generated edit script applying cell-level value changes
# derived from the unified diff of Golem_Profits.xlsx (8 sheets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3041.6667
$ws.Range("I40").Value = 3041.6667
$ws.Range("K40").Value = 3041.6667
$ws.Range("M40").Value = -2866.6667
$ws.Range("H51").Value = 92856.14
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("H88").Value = 950
$ws.Range("I88").Value = 900
$ws.Range("K88").Value = 900
$ws.Range("M88").Value = -494
$ws.Range("H91").Value = 950
$ws.Range("I91").Value = 900
$ws.Range("K91").Value = 900
$ws.Range("M91").Value = 504
$ws.Range("H112").Value = 2450
$ws.Range("I112").Value = 900
$ws.Range("K112").Value = 2700
$ws.Range("M112").Value = -1592
$ws.Range("H132").Value = 2575.75
$ws.Range("I132").Value = 768
$ws.Range("K132").Value = 2304
$ws.Range("M132").Value = 226
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 1750
$ws.Range("I21").Value = 1750
$ws.Range("K21").Value = 1750
$ws.Range("M21").Value = -1376
$ws.Range("H61").Value = 1674.2727
$ws.Range("I61").Value = 1694.2222
$ws.Range("J61").Value = 1584.5
$ws.Range("K61").Value = 1694.2222
$ws.Range("L61").Value = 1584.5
$ws.Range("M61").Value = -1482.2222
$ws.Range("N61").Value = -2008.5
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("H96").Value = 35868.8
$ws.Range("I96").Value = 26000
$ws.Range("J96").Value = 38336
$ws.Range("K96").Value = 26000
$ws.Range("L96").Value = 38336
$ws.Range("M96").Value = -23254
$ws.Range("N96").Value = -43828
$ws.Range("H133").Value = 100000
$ws.Range("I133").Value = 100000
$ws.Range("K133").Value = 100000
$ws.Range("M133").Value = -97470
$ws.Range("H136").Value = 1674.2727
$ws.Range("I136").Value = 1694.2222
$ws.Range("J136").Value = 1584.5
$ws.Range("K136").Value = 5082.6666
$ws.Range("L136").Value = 4753.5
$ws.Range("M136").Value = -2532.6666
$ws.Range("N136").Value = -9853.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 512
$ws.Range("I22").Value = 540
$ws.Range("J22").Value = 400
$ws.Range("K22").Value = 540
$ws.Range("L22").Value = 400
$ws.Range("M22").Value = -367
$ws.Range("N22").Value = -746
$ws.Range("H94").Value = 2394.1428
$ws.Range("I94").Value = 2001
$ws.Range("J94").Value = 2918.3333
$ws.Range("K94").Value = 2001
$ws.Range("L94").Value = 2918.3333
$ws.Range("M94").Value = -1550
$ws.Range("N94").Value = -3820.3333
$ws.Range("H105").Value = 2001
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 2001
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 2001
$ws.Range("M105").ClearContents()
$ws.Range("N105").Value = -5495
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 4994.5
$ws.Range("J4").Value = 4989
$ws.Range("L4").Value = 4989
$ws.Range("N4").Value = -5213
$ws.Range("H31").Value = 18935.875
$ws.Range("I31").Value = 7699.8
$ws.Range("J31").Value = 37662.668
$ws.Range("K31").Value = 7699.8
$ws.Range("L31").Value = 37662.668
$ws.Range("M31").Value = -7404.8
$ws.Range("N31").Value = -38252.668
$ws.Range("H34").Value = 18935.875
$ws.Range("I34").Value = 7699.8
$ws.Range("J34").Value = 37662.668
$ws.Range("K34").Value = 7699.8
$ws.Range("L34").Value = 37662.668
$ws.Range("M34").Value = -7497.8
$ws.Range("N34").Value = -38066.668
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1609.3541
$ws.Range("I4").Value = 1004.9667
$ws.Range("J4").Value = 2616.6667
$ws.Range("K4").Value = 3014.9001
$ws.Range("L4").Value = 7850.000100000001
$ws.Range("M4").Value = -2902.9001
$ws.Range("N4").Value = -8074.000100000001
$ws.Range("H17").Value = 245.70589
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("H39").Value = 5250
$ws.Range("I39").Value = 500
$ws.Range("J39").Value = 10000
$ws.Range("K39").Value = 1500
$ws.Range("L39").Value = 30000
$ws.Range("M39").Value = -1206
$ws.Range("N39").Value = -30588
$ws.Range("H55").Value = 1033.5454
$ws.Range("I55").Value = 816.9
$ws.Range("J55").Value = 3200
$ws.Range("K55").Value = 2450.7
$ws.Range("L55").Value = 9600
$ws.Range("M55").Value = -2273.7
$ws.Range("N55").Value = -9954
$ws.Range("H107").Value = 398.33334
$ws.Range("I107").Value = 398.33334
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1195.00002
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 724.9999800000001
$ws.Range("N107").ClearContents()
$ws.Range("H113").Value = 234.66667
$ws.Range("I113").Value = 315
$ws.Range("J113").Value = 194.5
$ws.Range("K113").Value = 945
$ws.Range("L113").Value = 583.5
$ws.Range("M113").Value = 1225
$ws.Range("N113").Value = -4923.5
$ws.Range("H134").Value = 942
$ws.Range("I134").Value = 942
$ws.Range("K134").Value = 2826
$ws.Range("M134").Value = 2244
$ws.Range("H137").Value = 6279.6
$ws.Range("J137").Value = 4724.5
$ws.Range("L137").Value = 14173.5
$ws.Range("N137").Value = -24373.5
$ws.Range("H138").Value = 1251.8
$ws.Range("I138").Value = 1251.8
$ws.Range("K138").Value = 3755.4
$ws.Range("M138").Value = 1384.6
$ws.Range("H141").Value = 2000
$ws.Range("I141").Value = 2000
$ws.Range("K141").Value = 6000
$ws.Range("M141").Value = -820
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 6142928.5
$ws.Range("J3").Value = 5000125
$ws.Range("L3").Value = 5000125
$ws.Range("N3").Value = -5000357
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("H93").Value = 71666.336
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 71666.336
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 71666.336
$ws.Range("M93").ClearContents()
$ws.Range("N93").Value = -75410.336
$ws.Range("H99").Value = 20140
$ws.Range("I99").Value = 12500
$ws.Range("J99").Value = 35420
$ws.Range("K99").Value = 12500
$ws.Range("L99").Value = 35420
$ws.Range("M99").Value = -10254
$ws.Range("N99").Value = -39912
$ws.Range("H102").Value = 45684.11
$ws.Range("I102").Value = 57909.715
$ws.Range("J102").Value = 2894.5
$ws.Range("K102").Value = 57909.715
$ws.Range("L102").Value = 2894.5
$ws.Range("M102").Value = -56287.715
$ws.Range("N102").Value = -6138.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 405
$ws.Range("I18").Value = 405
$ws.Range("K18").Value = 405
$ws.Range("M18").Value = -233
$ws.Range("H22").Value = 3095.6
$ws.Range("J22").Value = 3095.6
$ws.Range("L22").Value = 3095.6
$ws.Range("N22").Value = -3685.6
$ws.Range("H27").Value = 3095.6
$ws.Range("J27").Value = 3095.6
$ws.Range("L27").Value = 3095.6
$ws.Range("N27").Value = -3309.6
$ws.Range("H31").Value = 3240.3333
$ws.Range("J31").Value = 3977.2727
$ws.Range("L31").Value = 3977.2727
$ws.Range("N31").Value = -4473.2727
$ws.Range("H40").Value = 26823.8
$ws.Range("I40").Value = 22677.715
$ws.Range("K40").Value = 22677.715
$ws.Range("M40").Value = -22541.715
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 37000000
$ws.Range("I3").Value = 37000000
$ws.Range("K3").Value = 37000000
$ws.Range("M3").Value = -36999886
$ws.Range("H26").Value = 1000
$ws.Range("I26").Value = 1000
$ws.Range("K26").Value = 1000
$ws.Range("M26").Value = -707
$ws.Range("H132").Value = 1201.3334
$ws.Range("I132").Value = 802
$ws.Range("K132").Value = 2406
$ws.Range("M132").Value = 124
$ws.Range("H136").Value = 9372.583000000001
$ws.Range("I136").Value = 8696.4
$ws.Range("J136").Value = 9855.571
$ws.Range("K136").Value = 26089.2
$ws.Range("L136").Value = 29566.713
$ws.Range("M136").Value = -23539.2
$ws.Range("N136").Value = -34666.713
